$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C8").Value = 1.5
$ws.Range("D8").Value = "Autoupdating agenda list."
$ws.Range("D9").Select() | Out-Null
